# Generate Report for Handoff
# Updates the localization-status report so the Overview sheet reflects the
# latest handoff run: status text and timestamps move forward, and the
# status columns widen to fit the new "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status cells + the handoff generation date.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-07 01:18:12"

# zh-cn sheet: status + handoff datetime for the just-generated package.
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-07 01:18:03"

# de-de sheet: status (shares the same shared-string slot as zh-cn's).
$dede.Range("C2").Value = "Ready for handoff"

# Widen the status columns so the longer "Ready for handoff" text fits.
# (ColumnWidth is quantized to the host's pixel grid on write, same as real
# Excel; 16.3333... is the input that lands on the grid point nearest the
# target ~17.22-character width.)
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
